$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the extra repeated "value" header cells (C1:F1)
$ws.Range("C1:F1").ClearContents()

# Row 8: "Model" -> "production_function" (same "Sigmoid" value in B8)
$ws.Cells.Item(8,1).Value = "production_function"
$ws.Cells.Item(8,1).Style = $ws.Cells.Item(1,1).Style

# Insert a new row 9 for "L_curve" / 1
$ws.Rows("9:9").Insert()
$ws.Cells.Item(9,1).Value = "L_curve"
$ws.Cells.Item(9,1).Style = $ws.Cells.Item(1,1).Style
$ws.Cells.Item(9,2).Value = 1
$ws.Cells.Item(9,2).Style = $ws.Cells.Item(2,2).Style

# Remove the old "Deletion" row (now shifted down to row 17 because of the insert above)
$ws.Rows("17:17").Delete()

# Make this sheet the active / selected one, with C9 as the active cell
$ws.Activate()
$ws.Range("C9").Select()
